# Auto-generated Excel COM-interop script applying the Ravana_Profits update.
# For every affected (sheet, row) the H:N "price/profit" columns are refreshed
# to the latest Universalis snapshot values pulled in by the scheduled runner.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H113").Value = 3221.2222
$ws.Range("I113").Value = 3123.875
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 3123.875
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = 130.125
$ws.Range("N113").Value = -10508

$ws.Range("H116").Value = 4408.8
$ws.Range("J116").Value = 3523.5
$ws.Range("L116").Value = 3523.5
$ws.Range("N116").Value = -10407.5

$ws.Range("H132").Value = 1202.8948
$ws.Range("I132").Value = 1254
$ws.Range("K132").Value = 3762
$ws.Range("M132").Value = -1232

$ws.Range("H137").Value = 4026.6
$ws.Range("J137").Value = 6125.909
$ws.Range("L137").Value = 18377.727
$ws.Range("N137").Value = -23477.727

$ws.Range("H138").Value = 5676.4688
$ws.Range("J138").Value = 9922.4375
$ws.Range("L138").Value = 29767.3125
$ws.Range("N138").Value = -40047.3125


# ---- Sheet: ARM ----
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 3360.5264
$ws.Range("I32").Value = 2991.743
$ws.Range("K32").Value = 2991.743
$ws.Range("M32").Value = -2704.743

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0


# ---- Sheet: BSM ----
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H20").Value = 2925.5
$ws.Range("I20").Value = 3384.5
$ws.Range("J20").Value = 2237
$ws.Range("K20").Value = 3384.5
$ws.Range("L20").Value = 2237
$ws.Range("M20").Value = -3137.5
$ws.Range("N20").Value = -2731

$ws.Range("H82").Value = 37822.832
$ws.Range("I82").Value = 12652.333
$ws.Range("J82").Value = 62993.332
$ws.Range("K82").Value = 12652.333
$ws.Range("L82").Value = 62993.332
$ws.Range("M82").Value = -12269.333
$ws.Range("N82").Value = -63759.332

$ws.Range("H85").Value = 37822.832
$ws.Range("I85").Value = 12652.333
$ws.Range("J85").Value = 62993.332
$ws.Range("K85").Value = 12652.333
$ws.Range("L85").Value = 62993.332
$ws.Range("M85").Value = -11326.333
$ws.Range("N85").Value = -65645.33199999999

$ws.Range("H132").Value = 114499
$ws.Range("J132").Value = 114499
$ws.Range("L132").Value = 114499
$ws.Range("N132").Value = -124619


# ---- Sheet: CRP ----
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 1381.619
$ws.Range("I31").Value = 1486.5625
$ws.Range("J31").Value = 1045.8
$ws.Range("K31").Value = 1486.5625
$ws.Range("L31").Value = 1045.8
$ws.Range("M31").Value = -1191.5625
$ws.Range("N31").Value = -1635.8

$ws.Range("H34").Value = 1381.619
$ws.Range("I34").Value = 1486.5625
$ws.Range("J34").Value = 1045.8
$ws.Range("K34").Value = 1486.5625
$ws.Range("L34").Value = 1045.8
$ws.Range("M34").Value = -1284.5625
$ws.Range("N34").Value = -1449.8

$ws.Range("H58").Value = 1299
$ws.Range("I58").Value = 1299
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1299
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -1096

$ws.Range("H132").Value = 3413.182
$ws.Range("I132").Value = 3190.75
$ws.Range("J132").Value = 4658.8
$ws.Range("K132").Value = 9572.25
$ws.Range("L132").Value = 13976.4
$ws.Range("M132").Value = -7042.25
$ws.Range("N132").Value = -19036.4

$ws.Range("H134").Value = 4455.5835
$ws.Range("I134").Value = 4619
$ws.Range("J134").Value = 3965.3333
$ws.Range("K134").Value = 13857
$ws.Range("L134").Value = 11895.9999
$ws.Range("M134").Value = -11322
$ws.Range("N134").Value = -16965.9999

$ws.Range("H136").Value = 1299
$ws.Range("I136").Value = 1299
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3897
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -1347


# ---- Sheet: CUL ----
$ws = $wb.Sheets.Item("CUL")
$ws.Range("H7").Value = 155.44444
$ws.Range("I7").Value = 99.85714
$ws.Range("J7").Value = 350
$ws.Range("K7").Value = 299.57142
$ws.Range("L7").Value = 1050
$ws.Range("M7").Value = -187.57142
$ws.Range("N7").Value = -1274

$ws.Range("H118").Value = 4509.2646

$ws.Range("H138").Value = 1500
$ws.Range("I138").Value = 1500
$ws.Range("K138").Value = 4500
$ws.Range("M138").Value = 640


# ---- Sheet: GSM ----
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H70").Value = 7750
$ws.Range("I70").Value = 7750
$ws.Range("K70").Value = 7750
$ws.Range("M70").Value = -7480

$ws.Range("H73").Value = 7750
$ws.Range("I73").Value = 7750
$ws.Range("K73").Value = 7750
$ws.Range("M73").Value = -6814

$ws.Range("H80").Value = 4935.778
$ws.Range("I80").Value = 5004.25
$ws.Range("J80").Value = 4881
$ws.Range("K80").Value = 5004.25
$ws.Range("L80").Value = 4881
$ws.Range("M80").Value = -4006.25
$ws.Range("N80").Value = -6877

$ws.Range("H83").Value = 4935.778
$ws.Range("I83").Value = 5004.25
$ws.Range("J83").Value = 4881
$ws.Range("K83").Value = 25021.25
$ws.Range("L83").Value = 24405
$ws.Range("M83").Value = -20029.25
$ws.Range("N83").Value = -34389

$ws.Range("H97").Value = 1533
$ws.Range("I97").Value = 1533
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1533
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -1037


# ---- Sheet: LTW ----
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H16").Value = 799.6667
$ws.Range("I16").Value = 799.6667
$ws.Range("K16").Value = 799.6667
$ws.Range("M16").Value = -629.6667

$ws.Range("H82").Value = 1756
$ws.Range("I82").Value = 1769.8
$ws.Range("J82").Value = 1733
$ws.Range("K82").Value = 1769.8
$ws.Range("L82").Value = 1733
$ws.Range("M82").Value = -1408.8
$ws.Range("N82").Value = -2455

$ws.Range("H85").Value = 1756
$ws.Range("I85").Value = 1769.8
$ws.Range("J85").Value = 1733
$ws.Range("K85").Value = 1769.8
$ws.Range("L85").Value = 1733
$ws.Range("M85").Value = -521.8
$ws.Range("N85").Value = -4229

$ws.Range("H93").Value = 1397.3334
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H132").Value = 2709.476
$ws.Range("I132").Value = 1576.5
$ws.Range("K132").Value = 4729.5
$ws.Range("M132").Value = -2199.5

$ws.Range("H136").Value = 3797.5
$ws.Range("I136").Value = 3717.3333
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 11151.9999
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -8601.999899999999
$ws.Range("N136").Value = -20100


# ---- Sheet: WVR ----
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H62").Value = 16999.666
$ws.Range("I62").Value = 13500
$ws.Range("K62").Value = 13500
$ws.Range("M62").Value = -12876

$ws.Range("H65").Value = 16999.666
$ws.Range("I65").Value = 13500
$ws.Range("K65").Value = 67500
$ws.Range("M65").Value = -64380

$ws.Range("H81").Value = 5767.125
$ws.Range("I81").Value = 688.8333
$ws.Range("J81").Value = 21002
$ws.Range("K81").Value = 1377.6666
$ws.Range("L81").Value = 42004
$ws.Range("M81").Value = -316.6666
$ws.Range("N81").Value = -44126

$ws.Range("H84").Value = 5767.125
$ws.Range("I84").Value = 688.8333
$ws.Range("J84").Value = 21002
$ws.Range("K84").Value = 6888.333000000001
$ws.Range("L84").Value = 210020
$ws.Range("M84").Value = -1584.333000000001
$ws.Range("N84").Value = -220628

